# "Move gridcapacitygenerationlimit to specs"
#
# Adds a new parameter column, NewGridGenerationCapacityTimestepLimit, to
# the Malawi specs sheet: a new header in AJ1 (styled like the other
# bold/boxed headers, but with only a left/right border) and its value
# (9999999) in AJ2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + value in column AJ (the sheet currently ends at AI).
$ws.Range("AJ1").Value = "NewGridGenerationCapacityTimestepLimit"
$ws.Range("AJ2").Value = 9999999

# Match the look of the other header cells (bold, centered/top aligned)
# but give AJ1 a left/right-only box instead of the full four-side border
# used by the rest of row 1.
$ws.Range("AJ1").Font.Bold = $true
$ws.Range("AJ1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("AJ1").VerticalAlignment = -4160    # xlTop
$ws.Range("AJ1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft  = xlContinuous
$ws.Range("AJ1").Borders.Item(10).LineStyle = 1  # xlEdgeRight = xlContinuous

# Leave the new header cell selected/active, as in the authored workbook.
$ws.Range("AJ1").Select()
